# Weekly update: insert 2 new price rows for Femacal de La Calera - Frutilla
# This pushes the existing rows 119:239 down to 121:241 and fills the
# freed-up rows 119:120 with the newest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 119 (shifts 119-239 -> 121-241)
$ws.Rows("119:120").Insert()

# ---- New row 119 ----
$ws.Cells.Item(119, 1).Value  = 3
$ws.Cells.Item(119, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(119, 3).Value  = "Coquimbo"
$ws.Cells.Item(119, 4).Value  = 44587
$ws.Cells.Item(119, 5).Value  = 5
$ws.Cells.Item(119, 6).Value  = "Fruta"
$ws.Cells.Item(119, 7).Value  = 100101
$ws.Cells.Item(119, 8).Value  = "Berries"
$ws.Cells.Item(119, 9).Value  = 100112025
$ws.Cells.Item(119, 10).Value = "Frutilla"
$ws.Cells.Item(119, 11).Value = "Sin especificar"
$ws.Cells.Item(119, 12).Value = "Especial"
$ws.Cells.Item(119, 13).Value = 50
$ws.Cells.Item(119, 14).Value = 7000
$ws.Cells.Item(119, 15).Value = 7000
$ws.Cells.Item(119, 16).Value = 7000
$ws.Cells.Item(119, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(119, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(119, 19).Value = 1000
$ws.Cells.Item(119, 20).Value = 7

# ---- New row 120 ----
$ws.Cells.Item(120, 1).Value  = 3
$ws.Cells.Item(120, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(120, 3).Value  = "Coquimbo"
$ws.Cells.Item(120, 4).Value  = 44587
$ws.Cells.Item(120, 5).Value  = 5
$ws.Cells.Item(120, 6).Value  = "Fruta"
$ws.Cells.Item(120, 7).Value  = 100101
$ws.Cells.Item(120, 8).Value  = "Berries"
$ws.Cells.Item(120, 9).Value  = 100112025
$ws.Cells.Item(120, 10).Value = "Frutilla"
$ws.Cells.Item(120, 11).Value = "Sin especificar"
$ws.Cells.Item(120, 12).Value = "Primera"
$ws.Cells.Item(120, 13).Value = 180
$ws.Cells.Item(120, 14).Value = 5500
$ws.Cells.Item(120, 15).Value = 6000
$ws.Cells.Item(120, 16).Value = 5750
$ws.Cells.Item(120, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(120, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(120, 19).Value = 821
$ws.Cells.Item(120, 20).Value = 7

# Keep the date column formatted like the rest of column D
$ws.Range("D119:D120").NumberFormat = $ws.Range("D121").NumberFormat
